# "update In Class Demonstrations" - refresh the LF/FFR regression output
# table: new coefficient values for the LF Lag / FFR Lag rows, and drop the
# now-removed Constant / r2_adj rows entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-0.33***"
$ws.Range("C2").Value = "-0.03***"
$ws.Range("B3").Value = "7.684***"

# "0.071" would otherwise be auto-coerced to a numeric value by Excel, but
# the source data stores it as text - round-trip it through a scratch cell
# formatted as Text so it lands back in C3 as a string, then tidy the
# scratch cell back up.
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "0.071"
$ws.Range("Z1").Copy()
$ws.Range("C3").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("Z1").Clear()
$excel.CutCopyMode = 0

# The "Constant" and "r2_adj" rows (4 and 5) are no longer part of the
# summary table.
$ws.Range("A4:C5").EntireRow.Delete()
